$d = $word.ActiveDocument

# Rename the comment-delimiter tokens to bracket-delimiter tokens.
$d.Content.Find.Execute("<COMENTA>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<CORCHEA>", 2)
$d.Content.Find.Execute("<COMENTC>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<CORCHEC>", 2)
